$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.070.98"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "2.295.40"
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'113.01"
$ws.Range("E5").Value = "  -2.21%  "
$ws.Range("D6").Value = "'309.34"
$ws.Range("E6").Value = "  +2.54%  "
$ws.Range("D7").Value = "'0.635"
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").Value = "'44.55"
$ws.Range("E10").Value = "  -4.00%  "
$ws.Range("D11").Value = "'0.0930"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").Value = "'55.33"
$ws.Range("E12").Value = "  -2.71%  "
$ws.Range("D13").Value = "'8.85"
$ws.Range("E13").Value = "  -2.75%  "
$ws.Range("E14").Value = "  +21.13%  "
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").Value = "'15.56"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").Value = "2.642.35"
$ws.Range("E17").Value = "  +2.34%  "
$ws.Range("D18").Value = "2.298.41"
$ws.Range("E18").Value = "  +1.43%  "
$ws.Range("D19").Value = "43.058.07"
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").Value = "'7.23"
$ws.Range("E21").Value = "  -4.88%  "
$ws.Range("D22").Value = "'76.44"
$ws.Range("E22").Value = "  +3.00%  "
$ws.Range("D23").Value = "'3.56"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("E24").Value = "  +4.73%  "
$ws.Range("D25").Value = "'256.87"
$ws.Range("E25").Value = "  +10.33%  "
$ws.Range("E26").Value = "  -3.92%  "
$ws.Range("D27").Value = "'11.79"
$ws.Range("E27").Value = "  -3.51%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'39.27"
$ws.Range("E29").Value = "  -2.22%  "
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").Value = "'22.35"
$ws.Range("E31").Value = "  +4.98%  "
$ws.Range("D32").Value = "'173.94"
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("E33").Value = "  -2.95%  "
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("E36").Value = "  +7.03%  "
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("D38").Value = "'4.19"
$ws.Range("E38").Value = "  -8.35%  "
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("D42").Value = "'72.58"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("E44").Value = "  +6.54%  "
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").Value = "'12.50"
$ws.Range("E46").Value = "  -7.62%  "
$ws.Range("D47").Value = "'5.73"
$ws.Range("E47").Value = "  +2.37%  "
$ws.Range("D48").Value = "'109.01"
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("D49").Value = "'8.95"
$ws.Range("E49").Value = "  +3.91%  "
$ws.Range("E50").Value = "  -2.26%  "
$ws.Range("E51").Value = "  -0.51%  "
